# Generate Report for Handback
#
# - Marks all "Ready for handoff" status cells as handed back, in sync with en-US.
# - Records the handback timestamps for the zh-cn and de-de target files.
# - Adds "Latest Target File" / "Latest Handback File" hyperlink columns (F/G)
#   for both language sheets, mirroring the source-file / handoff-file links.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Helper: find the Address of the hyperlink anchored at a given A1 cell on a sheet.
function Get-HyperlinkAddress($sheet, $a1) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $a1) {
            return $h.Address
        }
    }
    return $null
}

# --- Overview sheet: refresh the zh-cn / de-de status columns -------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$zhA2 = Get-HyperlinkAddress $wsZh "`$A`$2"
$zhD2 = Get-HyperlinkAddress $wsZh "`$D`$2"
$zhA3 = Get-HyperlinkAddress $wsZh "`$A`$3"
$zhD3 = Get-HyperlinkAddress $wsZh "`$D`$3"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhA2, "", "", "4c752d07-9430-47c2-90b8-18c7a3dc4d5f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhD2, "", "", "4c752d07-9430-47c2-90b8-18c7a3dc4d5f.6397688340b0f8080e0db3b64037fd94e0aa0beb.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhA3, "", "", "a35c83ef-32bc-4ed9-99e5-805642bc0992.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhD3, "", "", "a35c83ef-32bc-4ed9-99e5-805642bc0992.36c4e7a7f3291a05a351340de378ec580bdbef88.zh-cn.xlf")

$wsZh.Range("H2").Value = "2016-03-17 03:46:25"
$wsZh.Range("H3").Value = "2016-03-17 03:46:25"

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$deA2 = Get-HyperlinkAddress $wsDe "`$A`$2"
$deD2 = Get-HyperlinkAddress $wsDe "`$D`$2"
$deA3 = Get-HyperlinkAddress $wsDe "`$A`$3"
$deD3 = Get-HyperlinkAddress $wsDe "`$D`$3"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deA2, "", "", "4c752d07-9430-47c2-90b8-18c7a3dc4d5f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deD2, "", "", "4c752d07-9430-47c2-90b8-18c7a3dc4d5f.6397688340b0f8080e0db3b64037fd94e0aa0beb.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deA3, "", "", "a35c83ef-32bc-4ed9-99e5-805642bc0992.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deD3, "", "", "a35c83ef-32bc-4ed9-99e5-805642bc0992.36c4e7a7f3291a05a351340de378ec580bdbef88.de-de.xlf")

$wsDe.Range("H2").Value = "2016-03-17 03:46:40"
$wsDe.Range("H3").Value = "2016-03-17 03:46:40"
